$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequence")

# Remove the "Notes Pending" note next to LSTM (B30).
$ws.Range("B30").ClearContents()

# Push the RNN-variant rows (Gated Recurrent Unit / Deep RNN /
# Bidirectional RNN) down by two rows (one slot is intentionally left
# blank at row 32).
$ws.Range("A33").Value = "Gated Recurrent Unit"
$ws.Range("A34").Value = "Deep RNN"
$ws.Range("A35").Value = "Bidirectional RNN"

# Append new Deep Learning topics being worked on.
$ws.Range("A36").Value = "LLM"
$ws.Range("A37").Value = "Encoder and Decoder"

# Insert a new "LSTM Architecture" entry right after LSTM.
$ws.Range("A31").Value = "LSTM Architecture"

$ws.Range("A38").Value = "Attention Mechanism"

$ws.Range("A32").ClearContents()

# Column A is now filled with shorter labels - shrink it to fit.
$ws.Columns.Item(1).ColumnWidth = 17.7

# Scroll the view down a bit and leave B38 selected.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B38").Select() | Out-Null
